$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "269.24"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.96"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.641"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.690"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8296"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01378"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1601"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08281"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03187"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09339"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.840"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001651"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04748"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006335"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005675"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001077"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.717"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.324"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3339"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04695"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007012"

$ws.Range("B42").Value = "CEJI"

$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003800"

$ws.Range("E42").Value = "41CEJICEJIWorstin24h"

$ws.Range("B43").Value = "BKEXToken"

$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1160"

$ws.Range("E43").Value = "42BKEXTokenBKK"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01191"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006266"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0009901"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9201"

$ws.Range("B49").Value = "BOLO"

$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002050"

$ws.Range("E49").Value = "48BOLOBOLO"

$ws.Range("B50").Value = "CryptobidCoin"

$ws.Range("C50").Value = "https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00001400"

$ws.Range("E50").Value = "49CryptobidCoinCBC"
